$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value() = $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.655.13"
$ws.Range("E2").Value() = "  -1.16%  "

Set-TextValue $ws.Range("D3") "1.881.18"
$ws.Range("E3").Value() = "  -0.47%  "

$ws.Range("E4").Value() = "  +0.08%  "

Set-TextValue $ws.Range("D5") "331.18"

$ws.Range("E6").Value() = "  +0.08%  "

Set-TextValue $ws.Range("D7") "0.4730"
$ws.Range("E7").Value() = "  +3.05%  "

Set-TextValue $ws.Range("D8") "0.3966"
$ws.Range("E8").Value() = "  +0.70%  "

Set-TextValue $ws.Range("D9") "47.15"
$ws.Range("E9").Value() = "  -4.11%  "

Set-TextValue $ws.Range("D10") "0.08036"
$ws.Range("E10").Value() = "  -2.49%  "

Set-TextValue $ws.Range("D11") "1.024"
$ws.Range("E11").Value() = "  -1.37%  "

Set-TextValue $ws.Range("D12") "21.77"
$ws.Range("E12").Value() = "  -0.14%  "

Set-TextValue $ws.Range("D13") "1.877.79"
$ws.Range("E13").Value() = "  -0.81%  "

Set-TextValue $ws.Range("D14") "5.968"
$ws.Range("E14").Value() = "  -0.27%  "

Set-TextValue $ws.Range("D15") "7.170"
$ws.Range("E15").Value() = "  -2.25%  "

$ws.Range("E16").Value() = "  +0.09%  "

Set-TextValue $ws.Range("D17") "86.92"
$ws.Range("E17").Value() = "  -2.65%  "

Set-TextValue $ws.Range("D18") "0.00001043"
$ws.Range("E18").Value() = "  -1.07%  "

Set-TextValue $ws.Range("D19") "0.06625"
$ws.Range("E19").Value() = "  +0.80%  "

Set-TextValue $ws.Range("D20") "17.18"
$ws.Range("E20").Value() = "  -1.93%  "

Set-TextValue $ws.Range("D21") "1.002"
$ws.Range("E21").Value() = "  +0.05%  "

Set-TextValue $ws.Range("D22") "27.669.58"
$ws.Range("E22").Value() = "  -1.06%  "

Set-TextValue $ws.Range("D23") "5.501"
$ws.Range("E23").Value() = "  -2.40%  "

Set-TextValue $ws.Range("D24") "10.98"
$ws.Range("E24").Value() = "  -0.95%  "

Set-TextValue $ws.Range("D25") "2.309"
$ws.Range("E25").Value() = "  +0.23%  "

Set-TextValue $ws.Range("D26") "2.083.43"
$ws.Range("E26").Value() = "  -1.61%  "

Set-TextValue $ws.Range("D27") "155.83"
$ws.Range("E27").Value() = "  +1.01%  "

Set-TextValue $ws.Range("D28") "20.23"
$ws.Range("E28").Value() = "  +1.56%  "

Set-TextValue $ws.Range("D29") "2.093"
$ws.Range("E29").Value() = "  -0.73%  "

Set-TextValue $ws.Range("D30") "5.564"
$ws.Range("E30").Value() = "  -2.56%  "

Set-TextValue $ws.Range("D31") "122.21"
$ws.Range("E31").Value() = "  -1.13%  "

Set-TextValue $ws.Range("D32") "0.9673"
$ws.Range("E32").Value() = "  +0.84%  "

Set-TextValue $ws.Range("D33") "0.09541"
$ws.Range("E33").Value() = "  -0.05%  "

Set-TextValue $ws.Range("D34") "1.462"
$ws.Range("E34").Value() = "  -0.84%  "

Set-TextValue $ws.Range("D35") "3.632"
$ws.Range("E35").Value() = "  +0.07%  "

$ws.Range("E36").Value() = "  -2.96%  "

Set-TextValue $ws.Range("D37") "0.06121"
$ws.Range("E37").Value() = "  +0.18%  "

Set-TextValue $ws.Range("D38") "0.02253"
$ws.Range("E38").Value() = "  -1.21%  "

Set-TextValue $ws.Range("D39") "1.233"
$ws.Range("E39").Value() = "  -1.85%  "

Set-TextValue $ws.Range("D40") "8.143"
$ws.Range("E40").Value() = "  -5.24%  "

Set-TextValue $ws.Range("D41") "0.6005"
$ws.Range("E41").Value() = "  -1.71%  "

$ws.Range("E42").Value() = "  +0.02%  "

Set-TextValue $ws.Range("D43") "0.1899"
$ws.Range("E43").Value() = "  +0.06%  "

Set-TextValue $ws.Range("D44") "10.29"
$ws.Range("E44").Value() = "  -4.33%  "

Set-TextValue $ws.Range("D45") "0.5714"
$ws.Range("E45").Value() = "  -1.62%  "

Set-TextValue $ws.Range("D46") "1.249"
$ws.Range("E46").Value() = "  -4.64%  "

Set-TextValue $ws.Range("D47") "12.25"
$ws.Range("E47").Value() = "  -3.90%  "

Set-TextValue $ws.Range("D48") "3.407"
$ws.Range("E48").Value() = "  -0.45%  "

Set-TextValue $ws.Range("D49") "1.932"
$ws.Range("E49").Value() = "  -3.07%  "

Set-TextValue $ws.Range("D50") "0.06819"
$ws.Range("E50").Value() = "  -1.17%  "

Set-TextValue $ws.Range("D51") "110.85"
$ws.Range("E51").Value() = "  +0.34%  "
